# fix bug with splitting strings unicode
# Append a third line to the wrapped alt-label text in C8, grow the row to
# fit the extra line, and move the selection to C9 (as if the user had
# just finished editing C8 and pressed Enter).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testreg4")

$ws.Range("C8").Value = "altLabel-multi-line1`naltLabel-multi-line2`nline3"

$ws.Rows.Item(8).RowHeight = 43.2

$ws.Range("C9").Select()
